# STD 3A - results.xlsx
# 11/15/2016 Added Streamreader examples in Chapter Methods.
#
# The underlying data edit: one duplicate roster row (Валентин Първанов,
# faculty no. 1601681087, grade 5) is removed, and the remaining 17 rows
# of the Table1 listobject are re-sorted ascending by column B
# (Факултетен номер / faculty number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate roster row. In the original sheet this is row 6
# (A6:C6 = "Валентин Първанов", 1601681087, 5) - a duplicate of the
# correctly-graded row further down (A12:C12, styled, grade 6). Deleting
# the row shifts everything below it up by one and the bound table
# (Table1) shrinks from A1:C20 to A1:C19 automatically.
$ws.Rows(6).Delete()

# Re-sort the table body by the "Факултетен номер" column (column B),
# ascending, keeping the header row in place.
$tbl = $ws.ListObjects.Item(1)
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("B1:B19"))
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# Match the resulting selection left behind in the saved file.
$ws.Range("A14:XFD14").Select()
